$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/coin/link columns to stay plain text so Excel does not
# auto-convert numeric-looking strings (e.g. "261.48") into numbers, and
# so the values keep their exact original text formatting.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value2 = "26.465.03"
$ws.Cells.Item(2, 5).Value2 = "  -0.72%  "
$ws.Cells.Item(3, 4).Value2 = "1.838.13"
$ws.Cells.Item(3, 5).Value2 = "  -0.86%  "
$ws.Cells.Item(4, 5).Value2 = "  -0.13%  "
$ws.Cells.Item(5, 4).Value2 = "261.48"
$ws.Cells.Item(5, 5).Value2 = "  -1.29%  "
$ws.Cells.Item(6, 4).Value2 = "1.000"
$ws.Cells.Item(6, 5).Value2 = "  -0.05%  "
$ws.Cells.Item(7, 4).Value2 = "0.5390"
$ws.Cells.Item(7, 5).Value2 = "  +2.12%  "
$ws.Cells.Item(8, 4).Value2 = "0.3020"
$ws.Cells.Item(8, 5).Value2 = "  -7.04%  "
$ws.Cells.Item(9, 5).Value2 = "  +0.97%  "
$ws.Cells.Item(10, 4).Value2 = "17.73"
$ws.Cells.Item(10, 5).Value2 = "  -6.69%  "
$ws.Cells.Item(11, 4).Value2 = "0.7385"
$ws.Cells.Item(11, 5).Value2 = "  -5.74%  "
$ws.Cells.Item(12, 4).Value2 = "1.835.63"
$ws.Cells.Item(12, 5).Value2 = "  -1.11%  "
$ws.Cells.Item(13, 4).Value2 = "0.07207"
$ws.Cells.Item(13, 5).Value2 = "  -7.61%  "
$ws.Cells.Item(14, 4).Value2 = "89.31"
$ws.Cells.Item(14, 5).Value2 = "  +0.75%  "
$ws.Cells.Item(15, 4).Value2 = "4.989"
$ws.Cells.Item(15, 5).Value2 = "  -0.88%  "
$ws.Cells.Item(16, 5).Value2 = "  -0.13%  "
$ws.Cells.Item(17, 4).Value2 = "13.84"
$ws.Cells.Item(17, 5).Value2 = "  -1.36%  "
$ws.Cells.Item(18, 4).Value2 = "1.000"
$ws.Cells.Item(18, 5).Value2 = "  -0.06%  "
$ws.Cells.Item(19, 4).Value2 = "0.000007872"
$ws.Cells.Item(19, 5).Value2 = "  -1.38%  "
$ws.Cells.Item(20, 4).Value2 = "26.490.53"
$ws.Cells.Item(20, 5).Value2 = "  -0.70%  "
$ws.Cells.Item(21, 4).Value2 = "2.081.04"
$ws.Cells.Item(21, 5).Value2 = "  -0.28%  "
$ws.Cells.Item(22, 4).Value2 = "4.591"
$ws.Cells.Item(22, 5).Value2 = "  -1.32%  "
$ws.Cells.Item(23, 4).Value2 = "5.971"
$ws.Cells.Item(23, 5).Value2 = "  -0.74%  "
$ws.Cells.Item(24, 4).Value2 = "9.221"
$ws.Cells.Item(24, 5).Value2 = "  -2.96%  "
$ws.Cells.Item(25, 4).Value2 = "142.71"
$ws.Cells.Item(25, 5).Value2 = "  -0.32%  "
$ws.Cells.Item(26, 5).Value2 = "  +0.71%  "
$ws.Cells.Item(27, 4).Value2 = "1.688"
$ws.Cells.Item(27, 5).Value2 = "  -0.44%  "
$ws.Cells.Item(28, 4).Value2 = "16.97"
$ws.Cells.Item(28, 5).Value2 = "  -0.37%  "
$ws.Cells.Item(29, 4).Value2 = "110.53"
$ws.Cells.Item(29, 5).Value2 = "  -1.44%  "
$ws.Cells.Item(30, 4).Value2 = "4.225"
$ws.Cells.Item(30, 5).Value2 = "  +0.62%  "
$ws.Cells.Item(31, 4).Value2 = "0.08828"
$ws.Cells.Item(31, 5).Value2 = "  +1.07%  "
$ws.Cells.Item(32, 4).Value2 = "4.027"
$ws.Cells.Item(32, 5).Value2 = "  -2.19%  "
$ws.Cells.Item(33, 4).Value2 = "0.04810"
$ws.Cells.Item(33, 5).Value2 = "  -0.84%  "
$ws.Cells.Item(34, 4).Value2 = "2.914"
$ws.Cells.Item(34, 5).Value2 = "  +1.20%  "
$ws.Cells.Item(35, 4).Value2 = "0.7281"
$ws.Cells.Item(35, 5).Value2 = "  +0.73%  "
$ws.Cells.Item(36, 4).Value2 = "1.131"
$ws.Cells.Item(36, 5).Value2 = "  -0.18%  "
$ws.Cells.Item(37, 5).Value2 = "  -0.74%  "
$ws.Cells.Item(38, 4).Value2 = "2.257"
$ws.Cells.Item(38, 5).Value2 = "  -0.70%  "
$ws.Cells.Item(39, 4).Value2 = "0.01711"
$ws.Cells.Item(39, 5).Value2 = "  -4.86%  "
$ws.Cells.Item(40, 4).Value2 = "0.4722"
$ws.Cells.Item(40, 5).Value2 = "  -3.30%  "
$ws.Cells.Item(41, 4).Value2 = "0.9026"
$ws.Cells.Item(41, 5).Value2 = "  -0.49%  "
$ws.Cells.Item(42, 4).Value2 = "107.71"
$ws.Cells.Item(42, 5).Value2 = "  -3.15%  "
$ws.Cells.Item(43, 4).Value2 = "5.885"
$ws.Cells.Item(43, 5).Value2 = "  -1.54%  "
$ws.Cells.Item(44, 5).Value2 = "  -0.08%  "
$ws.Cells.Item(45, 4).Value2 = "7.387"
$ws.Cells.Item(45, 5).Value2 = "  -4.08%  "
$ws.Cells.Item(46, 2).Value2 = "Algorand"
$ws.Cells.Item(46, 3).Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46, 4).Value2 = "0.1241"
$ws.Cells.Item(46, 5).Value2 = "  -0.01%  "
$ws.Cells.Item(47, 2).Value2 = "EnergySwap"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value2 = "8.960"
$ws.Cells.Item(47, 5).Value2 = "  -0.82%  "
$ws.Cells.Item(48, 4).Value2 = "0.4064"
$ws.Cells.Item(48, 5).Value2 = "  -3.61%  "
$ws.Cells.Item(49, 4).Value2 = "34.76"
$ws.Cells.Item(49, 5).Value2 = "  -1.10%  "
$ws.Cells.Item(50, 4).Value2 = "0.05773"
$ws.Cells.Item(50, 5).Value2 = "  -2.02%  "
$ws.Cells.Item(51, 4).Value2 = "0.8916"
$ws.Cells.Item(51, 5).Value2 = "  +0.20%  "

# Restore default (General) formatting now that the text values are set,
# so no visible number formatting change is left behind.
$textRange.ClearFormats()
